# Update KL Rahul (c)† match-by-match batting figures (runs, balls, fours, sixes)
# so that the per-row stats line up with the corrected match order.
# Values are entered with a leading apostrophe so they remain text-typed
# (matching the original t="str" cell storage) rather than becoming numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "'61"
$ws.Range("D2").Value = "'49"
$ws.Range("E2").Value = "'1"
$ws.Range("F2").Value = "'5"

# Row 4
$ws.Range("C4").Value = "'15"
$ws.Range("D4").Value = "'11"
$ws.Range("E4").Value = "'1"

# Row 5
$ws.Range("C5").Value = "'28"
$ws.Range("D5").Value = "'25"
$ws.Range("E5").Value = "'4"
$ws.Range("F5").Value = "'0"

# Row 6
$ws.Range("C6").Value = "'29"
$ws.Range("D6").Value = "'27"
$ws.Range("E6").Value = "'3"
$ws.Range("F6").Value = "'1"

# Row 7
$ws.Range("C7").Value = "'46"
$ws.Range("D7").Value = "'41"
$ws.Range("E7").Value = "'3"
$ws.Range("F7").Value = "'2"

# Row 8
$ws.Range("C8").Value = "'21"
$ws.Range("D8").Value = "'19"
$ws.Range("E8").Value = "'2"

# Row 9
$ws.Range("C9").Value = "'17"
$ws.Range("D9").Value = "'19"
$ws.Range("F9").Value = "'0"

# Row 11
$ws.Range("C11").Value = "'27"
$ws.Range("D11").Value = "'27"
$ws.Range("E11").Value = "'2"
$ws.Range("F11").Value = "'1"

# Row 12
$ws.Range("C12").Value = "'132"
$ws.Range("D12").Value = "'69"
$ws.Range("E12").Value = "'14"
$ws.Range("F12").Value = "'7"
